$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "71.185.80"
$ws.Range("E2").Value = "  +6.40%  "

# Row 3
$ws.Range("D3").Value = "3.806.55"
$ws.Range("E3").Value = "  +23.42%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "615.01"
$ws.Range("E5").Value = "  +7.73%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.62"
$ws.Range("E6").Value = "  +2.35%  "

# Row 7
$ws.Range("D7").Value = "3.819.33"
$ws.Range("E7").Value = "  +23.85%  "

# Row 8
$ws.Range("E8").Value = "  -0.10%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.547"
$ws.Range("E9").Value = "  +6.92%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.172"
$ws.Range("E10").Value = "  +13.35%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.41"
$ws.Range("E11").Value = "  +0.29%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.506"
$ws.Range("E12").Value = "  +8.56%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.06"
$ws.Range("E13").Value = "  +14.90%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000264"
$ws.Range("E14").Value = "  +9.83%  "

# Row 15
$ws.Range("D15").Value = "4.443.67"
$ws.Range("E15").Value = "  +23.51%  "

# Row 16
$ws.Range("D16").Value = "3.808.95"
$ws.Range("E16").Value = "  +23.63%  "

# Row 17
$ws.Range("D17").Value = "71.248.73"
$ws.Range("E17").Value = "  +6.57%  "

# Row 18
$ws.Range("E18").Value = "  +1.63%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.64"
$ws.Range("E19").Value = "  +9.25%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "525.52"
$ws.Range("E20").Value = "  +8.19%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.96"
$ws.Range("E21").Value = "  +2.90%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.53"
$ws.Range("E22").Value = "  +24.26%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.754"
$ws.Range("E23").Value = "  +10.48%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "89.00"
$ws.Range("E24").Value = "  +6.92%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.52"
$ws.Range("E25").Value = "  +13.09%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.68"
$ws.Range("E26").Value = "  +8.07%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.02"
$ws.Range("E27").Value = "  +8.27%  "

# Row 28
$ws.Range("E28").Value = "  +0.09%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000127"
$ws.Range("E29").Value = "  +35.24%  "

# Row 30
$ws.Range("E30").Value = "  +10.46%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.94"
$ws.Range("E31").Value = "  +14.18%  "

# Row 32
$ws.Range("E32").Value = "  +2.97%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "32.69"
$ws.Range("E33").Value = "  +17.05%  "

# Row 34
$ws.Range("E34").Value = "  +4.72%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.15%  "

# Row 36
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.24"
$ws.Range("E36").Value = "  +12.42%  "

# Row 37
$ws.Range("B37").Value = "Mantle"
$ws.Range("C37").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.06"
$ws.Range("E37").Value = "  +12.27%  "

# Row 38
$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.345"
$ws.Range("E38").Value = "  +11.54%  "

# Row 39
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.23"
$ws.Range("E39").Value = "  +11.55%  "

# Row 40
$ws.Range("E40").Value = "  +8.16%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "51.89"
$ws.Range("E41").Value = "  +6.03%  "

# Row 42
$ws.Range("E42").Value = "  +8.79%  "

# Row 43
$ws.Range("D43").Value = "3.165.96"
$ws.Range("E43").Value = "  +13.38%  "

# Row 44
$ws.Range("E44").Value = "  +16.85%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "44.35"
$ws.Range("E45").Value = "  -6.33%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.85"
$ws.Range("E46").Value = "  +5.23%  "

# Row 47
$ws.Range("E47").Value = "  +8.30%  "

# Row 48
$ws.Range("E48").Value = "  +9.96%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "142.09"
$ws.Range("E49").Value = "  +5.78%  "

# Row 51
$ws.Range("E51").Value = "  +8.45%  "
